$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2 = @{ E=3; G=18.82499266666667; H=56.474978; I=0.06886869772378311; J=0.0688686977237831; K=3; M=1.378475333333333; N=4.135426; O=0.05609715574531157; P=0.05609715574531156; Q=25.94978804118089; R=233.548092370628; S=0.003863338062187845; T=0.003863338062187844 }
  3 = @{ E=3; G=18.82499266666667; H=56.474978; I=0.06886869772378311; J=0.0688686977237831; K=3; M=13.16176133333333; N=39.485284; O=0.5356188518899525; P=0.5356188518899525; Q=247.7700605804169; R=2229.930545223752; S=0.0368873728059689; T=0.03688737280596889 }
  4 = @{ E=3; G=18.82499266666667; H=56.474978; I=0.06886869772378311; J=0.0688686977237831; K=3; M=0.829094; N=2.487282; O=0.03374004171190829; P=0.03374004171190828; Q=15.60768846997733; R=140.469196229796; S=0.002323632733845245; T=0.002323632733845245 }
  5 = @{ E=3; G=18.82499266666667; H=56.474978; I=0.06886869772378311; J=0.0688686977237831; K=3; M=9.203668; N=27.611004; O=0.3745439506528278; P=0.3745439506528276; Q=173.2589826064347; R=1559.330843457912; S=0.02579435412178113; T=0.02579435412178112 }
  6 = @{ E=3; G=121.8208923333333; H=365.462677; I=0.4456653109566078; J=0.4456653109566078; K=3; M=1.378475333333333; N=4.135426; O=0.05609715574531157; P=0.05609715574531156; Q=167.9270951661558; R=1511.343856495402; S=0.02500055635901554; T=0.02500055635901553 }
  7 = @{ E=3; G=121.8208923333333; H=365.462677; I=0.4456653109566078; J=0.4456653109566078; K=3; M=13.16176133333333; N=39.485284; O=0.5356188518899525; P=0.5356188518899525; Q=1603.37751030503; R=14430.39759274527; S=0.238706742181757; T=0.2387067421817569 }
  8 = @{ E=3; G=121.8208923333333; H=365.462677; I=0.4456653109566078; J=0.4456653109566078; K=3; M=0.829094; N=2.487282; O=0.03374004171190829; P=0.03374004171190828; Q=101.0009709082127; R=909.0087381739141; S=0.01503676618122653; T=0.01503676618122652 }
  9 = @{ E=3; G=121.8208923333333; H=365.462677; I=0.4456653109566078; J=0.4456653109566078; K=3; M=9.203668; N=27.611004; O=0.3745439506528278; P=0.3745439506528276; Q=1121.199048499745; R=10090.79143649771; S=0.1669212462346089; T=0.1669212462346088 }
  10 = @{ E=3; G=87.673585; H=263.020755; I=0.3207419907481189; J=0.3207419907481188; K=3; M=1.378475333333333; N=4.135426; O=0.05609715574531157; P=0.05609715574531156; Q=120.8558743074033; R=1087.70286876663; S=0.01799271340905851; T=0.0179927134090585 }
  11 = @{ E=3; G=87.673585; H=263.020755; I=0.3207419907481189; J=0.3207419907481188; K=3; M=13.16176133333333; N=39.485284; O=0.5356188518899525; P=0.5356188518899525; Q=1153.938801007713; R=10385.44920906942; S=0.1717954568374052; T=0.1717954568374052 }
  12 = @{ E=3; G=87.673585; H=263.020755; I=0.3207419907481189; J=0.3207419907481188; K=3; M=0.829094; N=2.487282; O=0.03374004171190829; P=0.03374004171190828; Q=72.68964328199; R=654.20678953791; S=0.01082184814660203; T=0.01082184814660203 }
  13 = @{ E=3; G=87.673585; H=263.020755; I=0.3207419907481189; J=0.3207419907481188; K=3; M=9.203668; N=27.611004; O=0.3745439506528278; P=0.3745439506528276; Q=806.9185687097801; R=7262.26711838802; S=0.1201319723550532; T=0.1201319723550531 }
  14 = @{ E=3; G=45.02666966666666; H=135.080009; I=0.1647240005714903; J=0.1647240005714903; K=3; M=1.378475333333333; N=4.135426; O=0.05609715574531157; P=0.05609715574531156; Q=62.06815347764822; R=558.6133812988339; S=0.009240547915049682; T=0.009240547915049678 }
  15 = @{ E=3; G=45.02666966666666; H=135.080009; I=0.1647240005714903; J=0.1647240005714903; K=3; M=13.16176133333333; N=39.485284; O=0.5356188518899525; P=0.5356188518899525; Q=592.6302797875062; R=5333.672518087556; S=0.0882292800648215; T=0.08822928006482149 }
  16 = @{ E=3; G=45.02666966666666; H=135.080009; I=0.1647240005714903; J=0.1647240005714903; K=3; M=0.829094; N=2.487282; O=0.03374004171190829; P=0.03374004171190828; Q=37.33134166061533; R=335.982074945538; S=0.005557794650234486; T=0.005557794650234485 }
  17 = @{ E=3; G=45.02666966666666; H=135.080009; I=0.1647240005714903; J=0.1647240005714903; K=3; M=9.203668; N=27.611004; O=0.3745439506528278; P=0.3745439506528276; Q=414.4105187576707; R=3729.694668819036; S=0.06169637794138463; T=0.0616963779413846 }
}

$colMap = @{ "E"=5; "G"=7; "H"=8; "I"=9; "J"=10; "K"=11; "M"=13; "N"=14; "O"=15; "P"=16; "Q"=17; "R"=18; "S"=19; "T"=20 }

foreach ($rowNum in $data.Keys) {
    $rowData = $data[$rowNum]
    foreach ($col in $rowData.Keys) {
        $colIdx = $colMap[$col]
        $ws.Cells.Item([int]$rowNum, $colIdx).Value = $rowData[$col]
    }
}

Write-Output "Applied Thbs1-Sdc1 NATMI update (Dr Hou advice)"